$d = $word.ActiveDocument

function Set-Run([int]$start, [int]$len, [string]$newText) {
    $rng = $d.Range($start, $start + $len)
    $rng.Text = $newText
}

# --- Edit 3: '... com icone “X” que represe...' -> '... com icone “lixeira” que represe...' ---
$r3 = $d.Content
$r3.Find.Execute("seleciona o botão com ícone “X” que represe") | Out-Null
$base3 = $r3.Start
# right-to-left: run D, run C, run B (run A 'seleciona...icone' is untouched)
Set-Run ($base3 + 30) 13 " que represe"
Set-Run ($base3 + 29) 1 "“lixeira”"
Set-Run ($base3 + 27) 2 " "

# --- Edit 2: 'Ator seleciona opcao de “Manter Racas”.' -> 'Ator seleciona opcao com icone “olho”.' ---
$r2 = $d.Content
$r2.Find.Execute("Ator seleciona opção de “Manter Raças”.") | Out-Null
$base2 = $r2.Start
# right-to-left: run C, run B, run A
Set-Run ($base2 + 37) 2 "."
Set-Run ($base2 + 25) 12 "com ícone “olho”"
Set-Run $base2 25 "Ator seleciona opção "

# --- Edit 1: 'Ator pressiona botao “Cadastrar Raca”.' -> 'Ator pressiona botao com icone “+”.' ---
$r1 = $d.Content
$r1.Find.Execute("Ator pressiona botão “Cadastrar Raça”.") | Out-Null
$base1 = $r1.Start
# right-to-left: run C, run B, run A
Set-Run ($base1 + 36) 2 "“+”."
Set-Run ($base1 + 22) 14 "ícone "
Set-Run $base1 22 "Ator pressiona botão com "

